$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at A, shifting existing columns (A->B, B->C, C->D, D->E, E->F)
$ws.Columns.Item(1).Insert()

# Header row
$ws.Range("A1").Value = "Category"
$ws.Range("B1").Value = "Type"
$ws.Range("C1").Value = "PMI_low"
$ws.Range("D1").Value = "PMI_high"
$ws.Range("E1").Value = "Yield_low"
$ws.Range("F1").Value = "Yield_high"

# Data rows
$data = @(
    @("Catalysis",     " Efficient Reaction",      2,   5,     0.95, 0.999),
    @("Catalysis",     "OK Reaction",              10,  30,    0.75, 0.9),
    @("Purification",  "Terrible Chromatography",  100, 10000, 0.2,  0.5),
    @("Purification",  "Great Chromatography",     50,  200,   0.7,  0.9),
    @("Misc",          "Mediocre Telescope",       50,  120,   0.65, 0.85),
    @("Misc",          "Other",                    10,  10000, 0.5,  1)
)

$r = 2
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $r++
}

# Apply percent style/format to Yield columns for all data rows
$ws.Range("E2:F7").Style = "Percent"
$ws.Range("E2:F7").NumberFormat = "0.0%"

# Column widths (carried over from the original bestFit widths of the shifted columns)
$ws.Columns.Item(2).Width = 127.494140625
$ws.Columns.Item(3).Width = 49.482421875
$ws.Columns.Item(4).Width = 52.4970703125
$ws.Columns.Item(5).Width = 54.732421875
$ws.Columns.Item(6).Width = 57.7470703125

# Selection
$ws.Range("C10").Select()
